$d = $word.ActiveDocument

# The document has two "${PROD}" placeholders; the one we must touch is the
# one embedded in "...that manufactures ${PROD}. The facility has an area..."
# Find that unique sentence first so we unambiguously land on the right
# paragraph/run (the other ${PROD} lives in the "Plant Principal Product:"
# line and must stay untouched).
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute('that manufactures ${PROD}.', $true, $false, $false, $false, $false, $false, 1, $false, '', 0)
if (-not $found) {
    throw 'Could not find the target sentence containing ${PROD}'
}

$paragraphRange = $find.Parent.Paragraphs(1).Range

$expectedParagraphText = "`tAssessment `${LE} was conducted on `${VDATE}, at a company located in `${LOC}, that manufactures `${PROD}. The facility has an area of `${AREA} ft2 for its manufacturing and office areas. It was reported during the assessment that production occurs `${PROH}. The office operates `${OFOH}. `r"
if ($paragraphRange.Text -ne $expectedParagraphText) {
    throw "Paragraph text did not match the expected introduction sentence; got: $($paragraphRange.Text)"
}

# InsertXML on this runtime replaces whole paragraphs touched by the range,
# so rebuild this paragraph's exact original OOXML (run-for-run, rsid-for-
# rsid) and splice in the new runs only around the "${PROD}" run - every
# other run in the paragraph is byte-for-byte unchanged.
$originalParagraphXml = '<w:p w14:paraId="18B3C8D7" w14:textId="2037655F" w:rsidR="008E53FA" w:rsidRDefault="00000000" w:rsidP="00A60200"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Assessment </w:t></w:r><w:r w:rsidR="00CC3535"><w:t>${LE}</w:t></w:r><w:r><w:t xml:space="preserve"> was conducted on </w:t></w:r><w:r w:rsidR="002B7FFC"><w:t>$</w:t></w:r><w:r w:rsidR="001A7D0F"><w:t>{</w:t></w:r><w:r w:rsidR="002B7FFC"><w:t>VDATE}</w:t></w:r><w:r><w:t xml:space="preserve">, at a company located in </w:t></w:r><w:r w:rsidR="00D8654A"><w:t>$</w:t></w:r><w:r w:rsidR="004E1E83"><w:t>{LOC}</w:t></w:r><w:r><w:t xml:space="preserve">, that manufactures </w:t></w:r><w:r w:rsidR="00A00294"><w:t>${PROD}</w:t></w:r><w:r><w:t xml:space="preserve">. The facility has an area of </w:t></w:r><w:r w:rsidR="00EF0A56"><w:t>${AREA}</w:t></w:r><w:r><w:t xml:space="preserve"> ft</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve"> for its manufacturing and office areas. It was reported during the assessment that production occurs </w:t></w:r><w:r w:rsidR="00EF0A56"><w:t>${PROH}</w:t></w:r><w:r><w:t>. The office operates</w:t></w:r><w:r w:rsidR="00EF0A56"><w:t xml:space="preserve"> ${OFOH}</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p>'

$oldRun = '<w:r w:rsidR="00A00294"><w:t>${PROD}</w:t></w:r>'
if ($originalParagraphXml -notlike "*$oldRun*") {
    throw 'Expected run for ${PROD} not found in the known paragraph XML'
}

# Word's spell checker flags "PRODlower" as a misspelling, which is why it
# ends up isolated in its own run, wrapped in proofErr spell-check markers,
# with the surrounding "${" and "}" kept in separate runs (same rsid as the
# original ${PROD} run, since it is the same edit/session).
$newRuns = '<w:r w:rsidR="00A00294"><w:t>${</w:t></w:r>' + `
           '<w:proofErr w:type="spellStart"/>' + `
           '<w:r w:rsidR="00A00294"><w:t>PRODlower</w:t></w:r>' + `
           '<w:proofErr w:type="spellEnd"/>' + `
           '<w:r w:rsidR="00A00294"><w:t>}</w:t></w:r>'

$newParagraphXml = $originalParagraphXml.Replace($oldRun, $newRuns)

$paragraphRange.InsertXML($newParagraphXml)

Write-Output 'Replaced ${PROD} run with split ${PRODlower} runs + proofErr markers'
